# The four records in rows 15-18 of the "Artfynd" sheet have been
# cyclically rotated by the source system:
#   new row 15 <- old row 18 (Koralltaggsvamp / Hericium coralloides)
#   new row 16 <- old row 17 (Grönpyrola / Pyrola chlorantha)
#   new row 17 <- old row 15 (Zontaggsvamp / Hydnellum concrescens)
#   new row 18 <- old row 16 (Orange taggsvamp / Hydnellum aurantiacum)
# Apply the resulting cell-value changes directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 15 (now holds what used to be row 18's record) ----
$ws.Range("A15").Value = 111837758
$ws.Range("B15").Value = 90187
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 2014
$ws.Range("F15").Value = "Koralltaggsvamp"
$ws.Range("G15").Value = "Hericium coralloides"
$ws.Range("H15").Value = "(Scop.:Fr.) Pers."
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = "6"
$ws.Range("P15").Value = "Brotorp, hyggeskant, Sm"
$ws.Range("Q15").Value = 575673.5681218
$ws.Range("R15").Value = 6404513.458820416
$ws.Range("AC15").Value = "På asplåga."

# ---- Row 16 (now holds what used to be row 17's record) ----
$ws.Range("A16").Value = 111837675
$ws.Range("B16").Value = 103288
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 221144
$ws.Range("F16").Value = "Grönpyrola"
$ws.Range("G16").Value = "Pyrola chlorantha"
$ws.Range("H16").Value = "Sw."
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "10"
$ws.Range("J16").Value = "plantor/tuvor"
# this record carries an (empty) Kön/L column value along with it
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = ""
$ws.Range("P16").Value = "Brotorp, Långsjön, Sm"
$ws.Range("Q16").Value = 575781.9606960951
$ws.Range("R16").Value = 6404546.96767282

# ---- Row 17 (now holds what used to be row 15's record) ----
$ws.Range("A17").Value = 111837705
$ws.Range("B17").Value = 90662
$ws.Range("E17").Value = 4363
$ws.Range("F17").Value = "Zontaggsvamp"
$ws.Range("G17").Value = "Hydnellum concrescens"
$ws.Range("H17").Value = "(Pers.) Banker"
$ws.Range("J17").Value = "fruktkroppar"
# the (empty) Kön/L column value no longer belongs to this record
$ws.Range("L17").ClearContents()
$ws.Range("Q17").Value = 575795.3141537429
$ws.Range("R17").Value = 6404518.948622406

# ---- Row 18 (now holds what used to be row 16's record) ----
$ws.Range("A18").Value = 111837741
$ws.Range("B18").Value = 90658
$ws.Range("E18").Value = 4361
$ws.Range("F18").Value = "Orange taggsvamp"
$ws.Range("G18").Value = "Hydnellum aurantiacum"
$ws.Range("H18").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I18").NumberFormat = "@"
$ws.Range("I18").Value = "15"
$ws.Range("Q18").Value = 575653.9215098171
$ws.Range("R18").Value = 6404506.688862759
# the public comment no longer belongs to this record
$ws.Range("AC18").ClearContents()
